$d = $word.ActiveDocument

function Split-RunAt($pos, $len) {
    # Toggling Bold on/off forces the run-coalescing logic to keep this
    # span as its own run instead of merging with identically formatted
    # neighbours. Must be done only AFTER all Find/Replace calls, since a
    # later Find/Replace re-coalesces identically formatted runs in the
    # paragraph it touches.
    $r = $d.Range($pos, $pos + $len)
    $r.Bold = $true
    $r.Bold = $false
}

# ---------------------------------------------------------------------
# Step 1: perform all the text substitutions first.
# ---------------------------------------------------------------------

# "from 2015 to  Early September 2022" -> "from 2015 to Early September 2022"
$d.Content.Find.Execute("from 2015 to  Early September 2022", $true, $false, $false, $false, $false, $true, 1, $false, "from 2015 to Early September 2022", 2)

# "5,181" -> "6,574"
$d.Content.Find.Execute("Within our data set of 5,181", $true, $false, $false, $false, $false, $true, 1, $false, "Within our data set of 6,574", 2)

# "2020" -> "2022"
$d.Content.Find.Execute("from 2015 to 2020 in the United States", $true, $false, $false, $false, $false, $true, 1, $false, "from 2015 to 2022 in the United States", 2)

# ---------------------------------------------------------------------
# Step 2: now split the runs to match the target run layout. All finds
# are done, so text positions are stable from here on.
# ---------------------------------------------------------------------

$t = $d.Content.Text

# Split "...from 2015 to " | "Early September 2022"
$needle1 = "from 2015 to Early September 2022"
$prefix1 = "from 2015 to "
$idx1 = $t.IndexOf($needle1)
$splitPos1 = $idx1 + $prefix1.Length
$runLen1 = ("Early September 2022").Length
Split-RunAt $splitPos1 $runLen1

# Split "Within our data set of " | "6,574"
$prefix2 = "Within our data set of "
$idx2 = $t.IndexOf($prefix2 + "6,574")
$splitPos2 = $idx2 + $prefix2.Length
$runLen2 = ("6,574").Length
Split-RunAt $splitPos2 $runLen2

# Split "...from 2015 to 202" | "2" | " in the United States, is there a
# correlation between the U.S. state of observation and "
$needle3 = "from 2015 to 2022 in the United States"
$prefix3 = "from 2015 to 202"
$idx3 = $t.IndexOf($needle3)
$splitPos3 = $idx3 + $prefix3.Length
Split-RunAt $splitPos3 1
$splitPos3b = $splitPos3 + 1
$tailText3 = " in the United States, is there a correlation between the U.S. state of observation and "
$tailLen3 = $tailText3.Length
Split-RunAt $splitPos3b $tailLen3
